$wb = $excel.ActiveWorkbook

# Sheet "OFF": update Row 3 ("R") target depth data
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 341
$wsOff.Range("C3").Value = 245
$wsOff.Range("D3").Value = 88
$wsOff.Range("E3").Value = 45
$wsOff.Range("G3").Value = 3

# Sheet "DEF": update Row 3 ("R") target depth data
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 338
$wsDef.Range("C3").Value = 229
$wsDef.Range("D3").Value = 67
$wsDef.Range("E3").Value = 33
